$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 248, shifting the existing rows
# (old 248..362) down to (250..364).
$ws.Range("A248:A249").EntireRow.Insert()

# Row 248: new weekly record (same dimensions/categorical values as the
# rest of the block, new date + volume/price figures).
$ws.Range("A248").Value = 9
$ws.Range("B248").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C248").Value = "Metropolitana"
$ws.Range("D248").Value = 44489
$ws.Range("E248").Value = 13
$ws.Range("F248").Value = 100112040
$ws.Range("G248").Value = "Cilantro"
$ws.Range("H248").Value = "Sin especificar"
$ws.Range("I248").Value = "Primera"
$ws.Range("J248").Value = 34
$ws.Range("K248").Value = 4000
$ws.Range("L248").Value = 4000
$ws.Range("M248").Value = 4000
$ws.Range("N248").Value = "$/caja 36 atados"
$ws.Range("O248").Value = "Región Metropolitana"
$ws.Range("P248").Value = 111
$ws.Range("Q248").Value = 36
$ws.Range("R248").Value = "Hortaliza"

# Row 249: second new weekly record.
$ws.Range("A249").Value = 9
$ws.Range("B249").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C249").Value = "Metropolitana"
$ws.Range("D249").Value = 44489
$ws.Range("E249").Value = 13
$ws.Range("F249").Value = 100112040
$ws.Range("G249").Value = "Cilantro"
$ws.Range("H249").Value = "Sin especificar"
$ws.Range("I249").Value = "Primera"
$ws.Range("J249").Value = 106
$ws.Range("K249").Value = 8000
$ws.Range("L249").Value = 10000
$ws.Range("M249").Value = 9000
$ws.Range("N249").Value = "$/docena de atados"
$ws.Range("O249").Value = "Región Metropolitana"
$ws.Range("P249").Value = 3000
$ws.Range("Q249").Value = 3
$ws.Range("R249").Value = "Hortaliza"
